# Requirements.xlsx update:
#  - Add a "Technologies used" mini-table in columns E:F (rows 1-13)
#  - Add bullet markers in column D for the new table rows
#  - Widen columns D & E
#  - Colour-code the existing "Completed?" column (green = completed, red = not completed)
#  - Colour-code the new "Used?" column (green = Yes, red = No)
#  - Make row 13 a touch taller (matches the new bold bullet there)
#  - Move the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "Technologies" / "Used?" table in columns E and F
# ---------------------------------------------------------------------------

$ws.Range("E1").Value = "Technologies"
$ws.Range("F1").Value = "Used?"
$ws.Range("E1").Font.Bold = $true
$ws.Range("F1").Font.Bold = $true

$technologies = @(
    @(2,  "Java 1.8",   "Yes"),
    @(3,  "Servlets",   "Yes"),
    @(4,  "JDBC",       "Yes"),
    @(5,  "SQL",        "Yes"),
    @(6,  "PL/SQL",     "Yes"),
    @(7,  "HTML/CSS",   "Yes"),
    @(8,  "Bootstrap",  "Yes"),
    @(9,  "JavaScript", "Yes"),
    @(10, "AJAX",       "No"),
    @(11, "Junit",      "Yes"),
    @(12, "Java Mail",  "Yes"),
    @(13, "Log4j",      "Yes")
)

# Write column E (technology names) first, then column F (Yes/No) - this
# keeps the shared-strings table in the same append order the original
# author ended up with (all technology names, then Yes, then No).
foreach ($row in $technologies) {
    $ws.Cells.Item($row[0], 5).Value = $row[1]
}

foreach ($row in $technologies) {
    $r = $row[0]
    $used = $row[2]
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = $used

    if ($used -eq "Yes") {
        $cell.Interior.Color = 5287936
    } else {
        $cell.Interior.Color = 255
    }
}

# ---------------------------------------------------------------------------
# 2) Bullet markers in column D next to the new table (rows 2-13)
# ---------------------------------------------------------------------------

$ws.Range("D2").HorizontalAlignment = -4131
$ws.Range("D2").VerticalAlignment = -4108
$ws.Range("D2").IndentLevel = 5

$ws.Range("D3:D12").Font.Name = "Symbol"
$ws.Range("D3:D12").Font.Size = 11
$ws.Range("D3:D12").HorizontalAlignment = -4131
$ws.Range("D3:D12").VerticalAlignment = -4108
$ws.Range("D3:D12").IndentLevel = 5

$ws.Range("D13").Font.Name = "Symbol"
$ws.Range("D13").Font.Size = 13
$ws.Range("D13").HorizontalAlignment = -4131
$ws.Range("D13").VerticalAlignment = -4108
$ws.Range("D13").IndentLevel = 5

# ---------------------------------------------------------------------------
# 3) Colour-code the existing "Completed?" column (B)
# ---------------------------------------------------------------------------

$completedCells = @(2,3,4,5,6,7,8,11,12,13,14,15,16,17,18,23)
foreach ($r in $completedCells) {
    $ws.Cells.Item($r, 2).Interior.Color = 5287936
}

$notCompletedCells = @(21,22,24,25)
foreach ($r in $notCompletedCells) {
    $ws.Cells.Item($r, 2).Interior.Color = 255
}

# ---------------------------------------------------------------------------
# 4) Column widths / row height
# ---------------------------------------------------------------------------

$ws.Columns.Item(4).ColumnWidth = 13.71
$ws.Columns.Item(5).ColumnWidth = 12.71
$ws.Rows.Item(13).RowHeight = 16.5

# ---------------------------------------------------------------------------
# 5) Selection (cosmetic, matches the saved workbook view)
# ---------------------------------------------------------------------------

$ws.Range("H11").Select()
